# This script applies the following change to the document:
#  1. In the paragraph about "cuantos dias tiene el mes", the text "...1,3,6,7,8,10 o 12..."
#     is changed to "...1,3,5,7,8,10 o 12..." (a single digit edit: 6 -> 5). Word records
#     this as the most recent edit location via the hidden "_GoBack" bookmark, placed
#     right after the freshly typed "5".
#  2. The old _GoBack bookmark (previously located right after the word "imprimir" in
#     the "invertir numero" exercise) is removed, since Word only ever keeps a single
#     _GoBack bookmark (tracking the most recent edit location).

$d = $word.ActiveDocument

# Make sure hidden bookmarks such as _GoBack are visible/addressable.
$d.Bookmarks.ShowHidden = $true

# --- Step 0: remove the (old) _GoBack bookmark first, while it is still unambiguous
#     (there is exactly one such bookmark in the original document, located right
#     after the word "imprimir" in the "invertir numero" exercise). Doing this before
#     inserting the new bookmark avoids ever having two same-named bookmarks at once. ---

if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
}

# --- Step 1: locate the run containing "...1,3,6,7,8,10 o 12..." and split it ---

$openQuote  = [char]0x201C
$closeQuote = [char]0x201D

$fullText = " un número e imprima para ese número cuántos días tiene el mes así: Si el número es el 11 o el 4 o el 6 o el 9, entonces imprima. " + `
            $openQuote + "Este mes tiene 30 días" + $closeQuote + ".  Si el número es igual a 2 entonces imprima " + `
            $openQuote + "el mes tiene 28 días" + $closeQuote + ",  si el mes es igual a 1,3,6,7,8,10 o 12, entonces imprima " + `
            $openQuote + "el mes tiene 31 días" + $closeQuote + ". Si el número no se encuentra entre 1 y 12 entonces imprima, " + `
            $openQuote + "lo siento el número no corresponde a un mes del año" + $closeQuote + "."

$probe = $d.Content
$found = $probe.Find.Execute($fullText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target paragraph text about days in a month."
}

# Use a *fresh* Range object built from the found boundaries (rather than reusing the
# Range that Find.Execute operated on) before handing it to InsertXML.
$targetRange = $d.Range($probe.Start, $probe.End)

$part1 = " un número e imprima para ese número cuántos días tiene el mes así: Si el número es el 11 o el 4 o el 6 o el 9, entonces imprima. " + `
         $openQuote + "Este mes tiene 30 días" + $closeQuote + ".  Si el número es igual a 2 entonces imprima " + `
         $openQuote + "el mes tiene 28 días" + $closeQuote + ",  si el mes es igual a 1,3,"
$part2 = "5"
$part3 = ",7,8,10 o 12, entonces imprima " + $openQuote + "el mes tiene 31 días" + $closeQuote + `
         ". Si el número no se encuentra entre 1 y 12 entonces imprima, " + $openQuote + `
         "lo siento el número no corresponde a un mes del año" + $closeQuote + "."

$rPr = "<w:rPr><w:rFonts w:eastAsia=`"Times New Roman`"/><w:szCs w:val=`"22`"/></w:rPr>"

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
       '<pkg:xmlData>' + `
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
       '<w:body><w:p>' + `
       '<w:r w:rsidRPr="00544487">' + $rPr + '<w:t xml:space="preserve">' + $part1 + '</w:t></w:r>' + `
       '<w:r>' + $rPr + '<w:t>' + $part2 + '</w:t></w:r>' + `
       '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
       '<w:r>' + $rPr + '<w:t>' + $part3 + '</w:t></w:r>' + `
       '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($xml)

# Sanity check: exactly one _GoBack bookmark should exist now, placed right after
# the newly inserted "5" (i.e. inside the run split we just performed).
if (-not $d.Bookmarks.Exists("_GoBack")) {
    throw "Expected a _GoBack bookmark after inserting the new run split, but none was found."
}
